$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting existing rows 41:133 down to 42:134
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new record
$ws.Cells.Item(41, 1).Value = 7
$ws.Cells.Item(41, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(41, 3).Value = "Ñuble"
$ws.Cells.Item(41, 4).Value = 44469
$ws.Cells.Item(41, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(41, 5).Value = 16
$ws.Cells.Item(41, 6).Value = 100112032
$ws.Cells.Item(41, 7).Value = "Zapallo italiano"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 60
$ws.Cells.Item(41, 11).Value = 17000
$ws.Cells.Item(41, 12).Value = 18000
$ws.Cells.Item(41, 13).Value = 17500
$ws.Cells.Item(41, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(41, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(41, 16).Value = 350
$ws.Cells.Item(41, 17).Value = 50
$ws.Cells.Item(41, 18).Value = "Hortaliza"
